$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2-13) got reshuffled: each row keeps its identifying
# columns (A,B,C,E-I,N,O,Q,R) but the date/volume/price columns
# (D,J,K,L,M,P) were rotated among the rows. Read all the original
# values first (using the Value() accessor, which actually invokes the
# getter), then write them back in the new order so source values
# aren't clobbered before being read.

$d2  = $ws.Range("D2").Value()
$j2  = $ws.Range("J2").Value()
$k2  = $ws.Range("K2").Value()
$l2  = $ws.Range("L2").Value()
$m2  = $ws.Range("M2").Value()
$p2  = $ws.Range("P2").Value()

$d3  = $ws.Range("D3").Value()
$j3  = $ws.Range("J3").Value()
$k3  = $ws.Range("K3").Value()
$l3  = $ws.Range("L3").Value()
$m3  = $ws.Range("M3").Value()
$p3  = $ws.Range("P3").Value()

$d4  = $ws.Range("D4").Value()
$j4  = $ws.Range("J4").Value()
$k4  = $ws.Range("K4").Value()
$l4  = $ws.Range("L4").Value()
$m4  = $ws.Range("M4").Value()
$p4  = $ws.Range("P4").Value()

$d5  = $ws.Range("D5").Value()
$j5  = $ws.Range("J5").Value()
$k5  = $ws.Range("K5").Value()
$l5  = $ws.Range("L5").Value()
$m5  = $ws.Range("M5").Value()
$p5  = $ws.Range("P5").Value()

$d7  = $ws.Range("D7").Value()
$j7  = $ws.Range("J7").Value()
$k7  = $ws.Range("K7").Value()
$l7  = $ws.Range("L7").Value()
$m7  = $ws.Range("M7").Value()
$p7  = $ws.Range("P7").Value()

$d8  = $ws.Range("D8").Value()
$j8  = $ws.Range("J8").Value()
$k8  = $ws.Range("K8").Value()
$l8  = $ws.Range("L8").Value()
$m8  = $ws.Range("M8").Value()
$p8  = $ws.Range("P8").Value()

$d9  = $ws.Range("D9").Value()
$j9  = $ws.Range("J9").Value()
$k9  = $ws.Range("K9").Value()
$l9  = $ws.Range("L9").Value()
$m9  = $ws.Range("M9").Value()
$p9  = $ws.Range("P9").Value()

$d10 = $ws.Range("D10").Value()
$j10 = $ws.Range("J10").Value()
$k10 = $ws.Range("K10").Value()
$l10 = $ws.Range("L10").Value()
$m10 = $ws.Range("M10").Value()
$p10 = $ws.Range("P10").Value()

$d11 = $ws.Range("D11").Value()
$j11 = $ws.Range("J11").Value()
$k11 = $ws.Range("K11").Value()
$l11 = $ws.Range("L11").Value()
$m11 = $ws.Range("M11").Value()
$p11 = $ws.Range("P11").Value()

$d13 = $ws.Range("D13").Value()
$j13 = $ws.Range("J13").Value()
$k13 = $ws.Range("K13").Value()
$l13 = $ws.Range("L13").Value()
$m13 = $ws.Range("M13").Value()
$p13 = $ws.Range("P13").Value()

# Row 2 <- old row 8
$ws.Range("D2").Value = $d8
$ws.Range("J2").Value = $j8
$ws.Range("K2").Value = $k8
$ws.Range("L2").Value = $l8
$ws.Range("M2").Value = $m8
$ws.Range("P2").Value = $p8

# Row 3 <- old row 9
$ws.Range("D3").Value = $d9
$ws.Range("J3").Value = $j9
$ws.Range("K3").Value = $k9
$ws.Range("L3").Value = $l9
$ws.Range("M3").Value = $m9
$ws.Range("P3").Value = $p9

# Row 4 <- old row 2
$ws.Range("D4").Value = $d2
$ws.Range("J4").Value = $j2
$ws.Range("K4").Value = $k2
$ws.Range("L4").Value = $l2
$ws.Range("M4").Value = $m2
$ws.Range("P4").Value = $p2

# Row 5 <- old row 11
$ws.Range("D5").Value = $d11
$ws.Range("J5").Value = $j11
$ws.Range("K5").Value = $k11
$ws.Range("L5").Value = $l11
$ws.Range("M5").Value = $m11
$ws.Range("P5").Value = $p11

# Row 6 unchanged

# Row 7 <- old row 10
$ws.Range("D7").Value = $d10
$ws.Range("J7").Value = $j10
$ws.Range("K7").Value = $k10
$ws.Range("L7").Value = $l10
$ws.Range("M7").Value = $m10
$ws.Range("P7").Value = $p10

# Row 8 <- old row 13
$ws.Range("D8").Value = $d13
$ws.Range("J8").Value = $j13
$ws.Range("K8").Value = $k13
$ws.Range("L8").Value = $l13
$ws.Range("M8").Value = $m13
$ws.Range("P8").Value = $p13

# Row 9 <- old row 5
$ws.Range("D9").Value = $d5
$ws.Range("J9").Value = $j5
$ws.Range("K9").Value = $k5
$ws.Range("L9").Value = $l5
$ws.Range("M9").Value = $m5
$ws.Range("P9").Value = $p5

# Row 10 <- old row 3
$ws.Range("D10").Value = $d3
$ws.Range("J10").Value = $j3
$ws.Range("K10").Value = $k3
$ws.Range("L10").Value = $l3
$ws.Range("M10").Value = $m3
$ws.Range("P10").Value = $p3

# Row 11 <- old row 7
$ws.Range("D11").Value = $d7
$ws.Range("J11").Value = $j7
$ws.Range("K11").Value = $k7
$ws.Range("L11").Value = $l7
$ws.Range("M11").Value = $m7
$ws.Range("P11").Value = $p7

# Row 12 unchanged

# Row 13 <- old row 4
$ws.Range("D13").Value = $d4
$ws.Range("J13").Value = $j4
$ws.Range("K13").Value = $k4
$ws.Range("L13").Value = $l4
$ws.Range("M13").Value = $m4
$ws.Range("P13").Value = $p4
